{"js": "// Update the \"Last Updated\" date: 06/10/2023 -> 06/26/2023\n// (the commit only forgot to update the day portion of the date).\nconst dateResults = context.document.body.search(\"06/10/2023\", { matchCase: true, matchWholeWord: false });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length === 0) {\n  throw new Error(\"Could not find the 'Last Updated' date (06/10/2023) to update.\");\n}\n\n// Only touch the day digits (\"10\" -> \"26\") within the matched date, so the\n// surrounding \"06/\", \"/\", \"2023\" text/runs are left untouched.\nconst dateRange = dateResults.items[0];\nconst dayResults = dateRange.search(\"10\", { matchCase: true, matchWholeWord: false });\ndayResults.load(\"items\");\nawait context.sync();\n\ndayResults.items[0].insertText(\"26\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the \"Last Updated\" date: 06/10/2023 -> 06/26/2023\n# (the commit only forgot to update the day portion of the date.)\n$d = $word.ActiveDocument\n\n# Locate the old date string so we only touch the day digits (\"10\" -> \"26\")\n# and leave the rest of \"Last Updated: 06/.../2023\" untouched.\n$dateRange = $d.Content\n$find = $dateRange.Find\n$find.ClearFormatting()\n$find.Text = \"06/10/2023\"\n$found = $find.Execute()\n\nif ($found) {\n    # \"06/10/2023\" -> the day digits \"10\" sit at offsets 3-4 within the match.\n    $dayStart = $dateRange.Start + 3\n    $dayRange = $d.Range($dayStart, $dayStart + 2)\n\n    $dayFind = $dayRange.Find\n    $dayFind.ClearFormatting()\n    $dayFind.Text = \"10\"\n    $dayFind.Replacement.ClearFormatting()\n    $dayFind.Replacement.Text = \"26\"\n    $dayFind.Execute([ref]$dayFind.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$dayFind.Replacement.Text, 2) | Out-Null\n} else {\n    Write-Output \"WARNING: date string '06/10/2023' not found\"\n}\n"}
